$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (dnn_n51_transpiled.qasm)
$ws.Range("B5").Value = 0.009443920072667891
$ws.Range("C5").Value = 0.004285694231455655
$ws.Range("D5").Value = 242

# Row 7 (sqrt18.qasm)
$ws.Range("B7").Value = [double]"5.658901510387511e-05"
$ws.Range("C7").Value = [double]"4.02702116170799e-05"
$ws.Range("D7").Value = 847
$ws.Range("E7").Value = 786

# Row 8 (dnn_n33_transpiled.qasm)
$ws.Range("B8").Value = 0.1021858996297044
$ws.Range("C8").Value = 0.07144515115937539
$ws.Range("D8").Value = 157
$ws.Range("E8").Value = 122

# Row 9 (qft_n18.qasm)
$ws.Range("B9").Value = 0.03547247095417511
$ws.Range("C9").Value = 0.03208437659567012
$ws.Range("D9").Value = 287
$ws.Range("E9").Value = 269

# Row 10 (DNN16.qasm)
$ws.Range("B10").Value = 0.4010494781974637
$ws.Range("C10").Value = 0.381939279255351
$ws.Range("D10").Value = 39
$ws.Range("E10").Value = 34

# Row 11 (QV_32.qasm)
$ws.Range("B11").Value = [double]"3.699410536551183e-09"
$ws.Range("C11").Value = [double]"2.330497946345709e-09"
$ws.Range("D11").Value = 1106
$ws.Range("E11").Value = 1059
$ws.Range("F11").Value = 1486
$ws.Range("G11").Value = 1479

# Row 14 (qaoa_n6_transpiled.qasm)
$ws.Range("C14").Value = 0.7603665634490904
$ws.Range("D14").Value = 31

# Row 15 (google_advantage.qasm)
$ws.Range("B15").Value = 0.7894068617535624
